$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row: num / RC / XRC / ERD
$ws2.Range("A1").Value = "num"
$ws2.Range("B1").Value = "RC"
$ws2.Range("C1").Value = "XRC"
$ws2.Range("D1").Value = "ERD"

$data = @(
    @(8, 1, 6.59, 6.59),
    @(16, 1, 7.05, 7.05),
    @(32, 1, 7.56, 7.56),
    @(64, 2, 8.65, 8.65),
    @(128, 3, 11.07, 11.07),
    @(256, 5, 14.25, 16.62),
    @(512, 7, 24.9, 26.8),
    @(1024, 10, 39.56, 46.97),
    @(2048, 17, 83.63, 93.22),
    @(4096, 28, 92.29, 95.42),
    @(8192, 29, 96.2, 97.81),
    @(16384, 33, 95.5, 96.88),
    @(32768, 38, 94.83, 95.99),
    @(65536, 31, 89.82, 91.06),
    @(131072, 29, 83.71, 84.56),
    @(262144, 37, 83.41, 83.88),
    @(524288, 30, 67.86, 68.25),
    @(1048576, 30, 67.77, 68.260000000000005)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws2.Activate()
$ws2.Range("D7").Select()
